$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 88; all rows 88..108 shift down to 89..109
$ws.Rows.Item(88).EntireRow.Insert()

$ws.Range("A88").Value = 11
$ws.Range("B88").Value = "Vega Monumental Concepción"
$ws.Range("C88").Value = "Bíobío"
$ws.Range("D88").Value = 44995
$ws.Range("D88").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E88").Value = 8
$ws.Range("F88").Value = 100112012
$ws.Range("G88").Value = "Espinaca"
$ws.Range("H88").Value = "Sin especificar"
$ws.Range("I88").Value = "Primera"
$ws.Range("J88").Value = 40
$ws.Range("K88").Value = 7500
$ws.Range("L88").Value = 8000
$ws.Range("M88").Value = 7750
$ws.Range("N88").Value = "`$/cuna 10 kilos"
$ws.Range("O88").Value = "Región Metropolitana"
$ws.Range("P88").Value = 775
$ws.Range("Q88").Value = 10
$ws.Range("R88").Value = "Hortaliza"
